$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 429.85
$ws.Range("I9").Value = 288.11765
$ws.Range("J9").Value = 1233
$ws.Range("K9").Value = 288.11765
$ws.Range("L9").Value = 1233
$ws.Range("M9").Value = -119.11765
$ws.Range("N9").Value = -1571
$ws.Range("H17").Value = 2229.1052
$ws.Range("J17").Value = 2286.2778
$ws.Range("L17").Value = 6858.8334
$ws.Range("N17").Value = -7194.8334
$ws.Range("H18").Value = 286
$ws.Range("I18").Value = 286
$ws.Range("K18").Value = 286
$ws.Range("M18").Value = -2
$ws.Range("H43").Value = 2579.9
$ws.Range("J43").Value = 2831.25
$ws.Range("L43").Value = 2831.25
$ws.Range("N43").Value = -2969.25
$ws.Range("H69").Value = 16693.957
$ws.Range("J69").Value = 17555.111
$ws.Range("L69").Value = 52665.333
$ws.Range("N69").Value = -54413.333
$ws.Range("H70").Value = 3164.8462
$ws.Range("I70").Value = 996.6667
$ws.Range("J70").Value = 3815.3
$ws.Range("K70").Value = 2990.0001
$ws.Range("L70").Value = 11445.9
$ws.Range("M70").Value = -2720.0001
$ws.Range("N70").Value = -11985.9
$ws.Range("H72").Value = 16693.957
$ws.Range("J72").Value = 17555.111
$ws.Range("L72").Value = 157995.999
$ws.Range("N72").Value = -166731.999
$ws.Range("H73").Value = 3164.8462
$ws.Range("I73").Value = 996.6667
$ws.Range("J73").Value = 3815.3
$ws.Range("K73").Value = 2990.0001
$ws.Range("L73").Value = 11445.9
$ws.Range("M73").Value = -2054.0001
$ws.Range("N73").Value = -13317.9
$ws.Range("H121").Value = 3826.2222
$ws.Range("J121").Value = 3826.2222
$ws.Range("L121").Value = 11478.6666
$ws.Range("N121").Value = -14972.6666
$ws.Range("H132").Value = 2419.415
$ws.Range("I132").Value = 1678.1428
$ws.Range("J132").Value = 11500
$ws.Range("K132").Value = 5034.428400000001
$ws.Range("L132").Value = 34500
$ws.Range("M132").Value = -2504.428400000001
$ws.Range("N132").Value = -39560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 532.7646999999999
$ws.Range("I2").Value = 526.14813
$ws.Range("K2").Value = 526.14813
$ws.Range("M2").Value = -413.14813
$ws.Range("H32").Value = 8558.041999999999
$ws.Range("I32").Value = 4688.4565
$ws.Range("K32").Value = 4688.4565
$ws.Range("M32").Value = -4401.4565
$ws.Range("H61").Value = 3528.261
$ws.Range("I61").Value = 3303.5881
$ws.Range("J61").Value = 4164.8335
$ws.Range("K61").Value = 3303.5881
$ws.Range("L61").Value = 4164.8335
$ws.Range("M61").Value = -3091.5881
$ws.Range("N61").Value = -4588.8335
$ws.Range("H74").Value = 6364.25
$ws.Range("I74").Value = 2916.9285
$ws.Range("K74").Value = 2916.9285
$ws.Range("M74").Value = -2042.9285
$ws.Range("H77").Value = 6364.25
$ws.Range("I77").Value = 2916.9285
$ws.Range("K77").Value = 14584.6425
$ws.Range("M77").Value = -10216.6425
$ws.Range("H116").Value = 532.7646999999999
$ws.Range("I116").Value = 526.14813
$ws.Range("K116").Value = 526.14813
$ws.Range("M116").Value = 1767.85187
$ws.Range("H122").Value = 2568.3704
$ws.Range("I122").Value = 2178.1738
$ws.Range("J122").Value = 4812
$ws.Range("K122").Value = 6534.5214
$ws.Range("L122").Value = 14436
$ws.Range("M122").Value = -4084.5214
$ws.Range("N122").Value = -19336
$ws.Range("H132").Value = 2220.3948
$ws.Range("I132").Value = 1778
$ws.Range("J132").Value = 3459.1
$ws.Range("K132").Value = 5334
$ws.Range("L132").Value = 10377.3
$ws.Range("M132").Value = -2804
$ws.Range("N132").Value = -15437.3
$ws.Range("H136").Value = 3528.261
$ws.Range("I136").Value = 3303.5881
$ws.Range("J136").Value = 4164.8335
$ws.Range("K136").Value = 9910.764299999999
$ws.Range("L136").Value = 12494.5005
$ws.Range("M136").Value = -7360.764299999999
$ws.Range("N136").Value = -17594.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 532.7646999999999
$ws.Range("I3").Value = 526.14813
$ws.Range("K3").Value = 526.14813
$ws.Range("M3").Value = -412.14813
$ws.Range("H86").Value = 1969.3513
$ws.Range("I86").Value = 1318.1538
$ws.Range("K86").Value = 1318.1538
$ws.Range("M86").Value = -195.1538
$ws.Range("H89").Value = 1969.3513
$ws.Range("I89").Value = 1318.1538
$ws.Range("K89").Value = 6590.769
$ws.Range("M89").Value = -974.7690000000002
$ws.Range("H105").Value = 1911858.6
$ws.Range("I105").Value = 2359473.5
$ws.Range("K105").Value = 2359473.5
$ws.Range("M105").Value = -2357726.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7398.1377
$ws.Range("I31").Value = 4697.512
$ws.Range("J31").Value = 13911.412
$ws.Range("K31").Value = 4697.512
$ws.Range("L31").Value = 13911.412
$ws.Range("M31").Value = -4402.512
$ws.Range("N31").Value = -14501.412
$ws.Range("H34").Value = 7398.1377
$ws.Range("I34").Value = 4697.512
$ws.Range("J34").Value = 13911.412
$ws.Range("K34").Value = 4697.512
$ws.Range("L34").Value = 13911.412
$ws.Range("M34").Value = -4495.512
$ws.Range("N34").Value = -14315.412
$ws.Range("H99").Value = 2274.9524
$ws.Range("I99").Value = 1853.3
$ws.Range("K99").Value = 1853.3
$ws.Range("M99").Value = -355.3
$ws.Range("H105").Value = 2131.4614
$ws.Range("I105").Value = 2320.818
$ws.Range("J105").Value = 1090
$ws.Range("K105").Value = 2320.818
$ws.Range("L105").Value = 1090
$ws.Range("M105").Value = -573.8180000000002
$ws.Range("N105").Value = -4584
$ws.Range("H122").Value = 1890.1111
$ws.Range("I122").Value = 1801.238
$ws.Range("J122").Value = 2014.5333
$ws.Range("K122").Value = 5403.714
$ws.Range("L122").Value = 6043.5999
$ws.Range("M122").Value = -2953.714
$ws.Range("N122").Value = -10943.5999
$ws.Range("H126").Value = 2274.9524
$ws.Range("I126").Value = 1853.3
$ws.Range("K126").Value = 5559.9
$ws.Range("M126").Value = -3089.9
$ws.Range("H132").Value = 4428.6904
$ws.Range("I132").Value = 2735.5293
$ws.Range("J132").Value = 11624.625
$ws.Range("K132").Value = 8206.5879
$ws.Range("L132").Value = 34873.875
$ws.Range("M132").Value = -5676.5879
$ws.Range("N132").Value = -39933.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 440.86957
$ws.Range("J92").Value = 447.72223
$ws.Range("L92").Value = 1343.16669
$ws.Range("N92").Value = -3839.16669
$ws.Range("H132").Value = 3074.5
$ws.Range("I132").Value = 3074.5
$ws.Range("K132").Value = 27670.5
$ws.Range("M132").Value = -25140.5
$ws.Range("H134").Value = 2131
$ws.Range("I134").Value = 2102.6155
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 6307.8465
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -1237.8465
$ws.Range("N134").Value = -17640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5031.8076
$ws.Range("J70").Value = 5037.08
$ws.Range("L70").Value = 5037.08
$ws.Range("N70").Value = -5577.08
$ws.Range("H73").Value = 5031.8076
$ws.Range("J73").Value = 5037.08
$ws.Range("L73").Value = 5037.08
$ws.Range("N73").Value = -6909.08
$ws.Range("H105").Value = 69949
$ws.Range("J105").Value = 69949
$ws.Range("L105").Value = 69949
$ws.Range("N105").Value = -76937
$ws.Range("H122").Value = 6325.75
$ws.Range("I122").Value = 3324
$ws.Range("J122").Value = 19333.334
$ws.Range("K122").Value = 9972
$ws.Range("L122").Value = 58000.00199999999
$ws.Range("M122").Value = -7522
$ws.Range("N122").Value = -62900.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5782.4
$ws.Range("I61").Value = 2184.3333
$ws.Range("J61").Value = 11179.5
$ws.Range("K61").Value = 2184.3333
$ws.Range("L61").Value = 11179.5
$ws.Range("M61").Value = -1982.3333
$ws.Range("N61").Value = -11583.5
$ws.Range("H82").Value = 1180.7693
$ws.Range("I82").Value = 776
$ws.Range("J82").Value = 1433.75
$ws.Range("K82").Value = 776
$ws.Range("L82").Value = 1433.75
$ws.Range("M82").Value = -415
$ws.Range("N82").Value = -2155.75
$ws.Range("H85").Value = 1180.7693
$ws.Range("I85").Value = 776
$ws.Range("J85").Value = 1433.75
$ws.Range("K85").Value = 776
$ws.Range("L85").Value = 1433.75
$ws.Range("M85").Value = 472
$ws.Range("N85").Value = -3929.75
$ws.Range("H113").Value = 5782.4
$ws.Range("I113").Value = 2184.3333
$ws.Range("J113").Value = 11179.5
$ws.Range("K113").Value = 2184.3333
$ws.Range("L113").Value = 11179.5
$ws.Range("M113").Value = -14.33329999999978
$ws.Range("N113").Value = -15519.5
$ws.Range("H132").Value = 6835.273
$ws.Range("I132").Value = 4711
$ws.Range("K132").Value = 14133
$ws.Range("M132").Value = -11603

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8335975
$ws.Range("I62").Value = 12502469
$ws.Range("J62").Value = 2987
$ws.Range("K62").Value = 12502469
$ws.Range("L62").Value = 2987
$ws.Range("M62").Value = -12501845
$ws.Range("N62").Value = -4235
$ws.Range("H65").Value = 8335975
$ws.Range("I65").Value = 12502469
$ws.Range("J65").Value = 2987
$ws.Range("K65").Value = 62512345
$ws.Range("L65").Value = 14935
$ws.Range("M65").Value = -62509225
$ws.Range("N65").Value = -21175
$ws.Range("H81").Value = 1511.2858
$ws.Range("I81").Value = 1359.909
$ws.Range("K81").Value = 2719.818
$ws.Range("M81").Value = -1658.818
$ws.Range("H84").Value = 1511.2858
$ws.Range("I84").Value = 1359.909
$ws.Range("K84").Value = 13599.09
$ws.Range("M84").Value = -8295.09
$ws.Range("H96").Value = 5718.9473
$ws.Range("J96").Value = 11999.429
$ws.Range("L96").Value = 11999.429
$ws.Range("N96").Value = -14745.429
$ws.Range("H132").Value = 3362.8086
$ws.Range("I132").Value = 2723.3777
$ws.Range("K132").Value = 8170.1331
$ws.Range("M132").Value = -5640.1331
